$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("D").Insert()
$ws.Columns("E").Copy()
$ws.Columns("D").PasteSpecial(-4122)
Write-Host "done"
